$d = $word.ActiveDocument

function Replace-Exact($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $true, $false, $false, $false,
                             $true, 1, $false, $newText, 2) | Out-Null
}

# Story 1: User logs in / can log into website and view stored information
Replace-Exact `
    "(5) User logs in using account to retrieve stored personal information. " `
    "(5) User can log into the website and view their stored information (i.e schedule and friends list). While not necessary for functionality this forms a base for the usability of the service and thus is high priority."

# Story 2: manual course input
Replace-Exact `
    "(5) User can input each course of their schedule manually." `
    "(5) User can manually input courses using a form and add the course to their schedule. This is vital to functionality."

# Story 3: automatic schedule input
Replace-Exact `
    "(3) Use can automatically input schedule into database without inputting individual courses." `
    "(3) User can automatically input their entire schedule by copying from student profile on the nest. This significantly improves usability but is not crucial to functionality."

# Story 4: view map of building
Replace-Exact `
    "(5) User can view a map of the building to aid in navigation" `
    "(5) User can view a map of the building to aid in navigation. This is a fundamental part of the system."

# Story 5: auto-generated path on map
Replace-Exact `
    "(2) Schedule automatically generates a path on map from class to class." `
    "(2) Schedule automatically generates a path on map for the entire user schedule. Very helpful and is the basis for the story 6, but not critical to functionality."

# Story 6: print map as pdf
Replace-Exact `
    "(1) User can print their auto generated map as pdf to save and use later." `
    "(1) User can print their auto generated map as pdf to save and use later. Reliant story 5, not crucial for functionality."

# Story 7: manual input two classes in text boxes
Replace-Exact `
    "(4) User can manually input two classes in text boxes and generate a path within one building" `
    "(4) User can manually input two classrooms into text boxes and generate a path within one building. Important to functionality."

# Story 8: interact with map / click on classrooms
Replace-Exact `
    "(1) User can interact with the map and click on classrooms to select for starting/ending navigation point." `
    "(1) User can interact with the map and click on classrooms to select for starting/ending navigation point. Enhancement to usability, but not necessary."

# Story 9: navigate between multiple buildings
Replace-Exact `
    "(1) User can navigate between multiple buildings on campus. (i.e. Classroom in Bellermine to classroom in Barbeline)" `
    "(1) User can navigate between multiple buildings on campus. (i.e. Classroom in Bellermine to classroom in Barbeline). Enhances usefulness of system as a whole and can easily be appended later."

# Story 10: printout of schedule
Replace-Exact `
    "(2) User can create a printout of their schedule in easy to view form." `
    "(2) User can create a printout of their schedule in easy to view form. Enhancement to usability, but not necessary."
